# Applies the workbook edit described in the commit:
# "Retirement Account Types. Fixes for 'How Much You Need' and 'How Long Does It Take'."
#
# Summary of changes:
#  - "Fixed Saving" sheet: C2:C37 (Growth = Start Balance * Rate) becomes one filled-down
#    (shared) formula instead of 37 separately authored formulas.
#  - "Raise Saving Yearly" sheet:
#      * Retirement budget input values in column F (rows 2-29) are all reduced by 1000
#        (and the dependent "Portfolio Needed" column G recalculates accordingly).
#      * C2:C29 (Growth) and G2:G29 (Portfolio Needed) become filled-down (shared) formulas.
#      * This sheet becomes the active/selected sheet & tab, with cell J5 selected.
#  - "Savings Rates" sheet: E2:E11 (Years To Retire) becomes one filled-down (shared) formula.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "Fixed Saving" sheet
# ---------------------------------------------------------------------------
$fixedSaving = $wb.Worksheets.Item("Fixed Saving")

# Re-enter the Growth formula across the whole column at once so it is saved as a single
# filled-down formula (C2:C37), matching B2:B37 * Rate for every row.
$fixedSaving.Range("C2:C37").Formula = "=B2*Rate"

# ---------------------------------------------------------------------------
# 2. "Raise Saving Yearly" sheet
# ---------------------------------------------------------------------------
$raiseSaving = $wb.Worksheets.Item("Raise Saving Yearly")

# Lower the planned retirement budget (column F) by 1000 for every populated row (2-29).
for ($r = 2; $r -le 29; $r++) {
    $budgetCell = $raiseSaving.Cells.Item($r, 6)   # column F
    $currentBudget = $budgetCell.Value()
    $budgetCell.Value = $currentBudget - 1000
}

# Re-enter the Growth and Portfolio Needed formulas across their full columns so they are
# saved as single filled-down formulas.
$raiseSaving.Range("C2:C29").Formula = "=B2*Rate"
$raiseSaving.Range("G2:G29").Formula = "=F2*YearsOfSpending"

# Make this the active sheet/tab, with J5 selected, as left by the author.
$raiseSaving.Activate()
$raiseSaving.Range("J5").Select()

# ---------------------------------------------------------------------------
# 3. "Savings Rates" sheet
# ---------------------------------------------------------------------------
$savingsRates = $wb.Worksheets.Item("Savings Rates")

# Re-enter the Years To Retire formula across the whole column at once so it is saved as a
# single filled-down formula (E2:E11).
$savingsRates.Range("E2:E11").Formula = "=LN(1+(Rate)*(YearsOfSpending*C2/A2))/LN(1+Rate)"
